# P3DFps.xlsx - add new FPS test rows (per-polygon AABB tests / texture
# quality fix runs) below the existing "Consolidate spans" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 37-39: continuation of the existing "Consolidate spans" shared-formula
# block (column A already carries style index 3 on these two pre-existing
# blank rows; row 39 is brand new so the number format is (re)applied
# explicitly to pick up the same style).
# NOTE: labels are registered into the shared-string table in the same
# first-use order as the source workbook ("Skip consolidate" before
# "Skip consolidate 4"), so row 38 is populated ahead of row 37.
$rows1 = @(
    @{ Row = 38; Label = "Skip consolidate";   Value = 166 },
    @{ Row = 37; Label = "Skip consolidate 4"; Value = 171 },
    @{ Row = 39; Label = "120*80";             Value = 371 }
)
foreach ($r in $rows1) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Label
    $ws.Range("A$row").NumberFormat = "#,##0"
    $ws.Range("B$row").Value = $r.Value
    $ws.Range("C$row").Formula = "=B$row/30"
    $ws.Range("D$row").Formula = '=B' + $row + '/$B$2'
}

# Remaining groups: each label gets its own pair of rows, separated by a
# blank spacer row (40, 43, 46, 49, 52 stay empty, matching the source).
$groups = @(
    @{ Label = "BB Test polys";   Rows = @(41, 42); Values = @(342, 186) },
    @{ Label = "BB Test polys 2"; Rows = @(44, 45); Values = @(347, 191) },
    @{ Label = "8x unroll";       Rows = @(47, 48); Values = @(338, 177) },
    @{ Label = "No render stats"; Rows = @(50, 51); Values = @(355, 235) },
    @{ Label = "Division factor"; Rows = @(53, 54); Values = @(370, 235) }
)

foreach ($g in $groups) {
    for ($i = 0; $i -lt $g.Rows.Count; $i++) {
        $row = $g.Rows[$i]
        $val = $g.Values[$i]
        $ws.Range("A$row").Value = $g.Label
        $ws.Range("A$row").NumberFormat = "#,##0"
        $ws.Range("B$row").Value = $val
        $ws.Range("C$row").Formula = "=B$row/30"
        $ws.Range("D$row").Formula = '=B' + $row + '/$B$2'
    }
}

# Match the author's final selection/viewport state.
$ws.Range("B54").Select()
